$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns K:O ---
$ws.Range("K1").Value = "maskEnabled"
$ws.Range("L1").Value = "a"
$ws.Range("M1").Value = "b"
$ws.Range("N1").Value = "c"
$ws.Range("O1").Value = "d"

# --- Row 2 (Georgia / training_lexical) ---
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("K2").Value = "yes"
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 20
$ws.Range("O2").Value = 20

# --- Row 3 (RobotoFlex / lexical_wo_driving_roboto) ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("K3").Value = "yes"
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 30
$ws.Range("O3").Value = 30

# --- Row 4 (Neue Frutiger World / lexical_wo_driving_neuefrutigerworld) ---
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("K4").Value = "yes"
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 30
$ws.Range("O4").Value = 30

# --- Row 5 (Eurostile / lexical_wo_driving_eurostile) ---
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("J5").Value = 2.14
$ws.Range("K5").Value = "yes"
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 30
$ws.Range("O5").Value = 30

# --- Row 6 (Georgia / full_task_training) ---
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("K6").Value = "no"
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 8

# --- Row 7 (RobotoFlex / full_task_roboto) ---
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("K7").Value = "no"
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 12
$ws.Range("O7").Value = 12

# --- Row 8 (Neue Frutiger World / full_task_neuefrutigerworld) ---
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("K8").Value = "no"
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 12
$ws.Range("O8").Value = 12

# --- Row 9 (Eurostile / full_task_eurostile) ---
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("J9").Value = 2.14
$ws.Range("K9").Value = "no"
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 10
$ws.Range("N9").Value = 12
$ws.Range("O9").Value = 12

# --- Selection moves to J17 in the saved file ---
$ws.Range("J17").Select() | Out-Null
